# Apply the "Plan" worksheet update:
#   - B2: "Complete Input Module" -> "Working On AI Module"
#   - B3: "Complete GPS Module"   -> "Complete Input Module"
#   - C5: (empty)                 -> "Complete GPS Module"
#   - Selection moves from D3 to B7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Working On AI Module"
$ws.Range("B3").Value = "Complete Input Module"
$ws.Range("C5").Value = "Complete GPS Module"

$ws.Range("B7").Select()
